$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing-cell fixes ---
$ws.Range("Q55").Value = 0
$ws.Range("Q57").Value = 0
$ws.Range("R354").Value = 0
$ws.Range("R355").Value = 0

# --- New rows 356:361 (weekly OHLCV bars for Jul-Aug 2024) ---
# row 356
$ws.Cells.Item(356, 1).Value = 45474
$ws.Cells.Item(356, 2).Value = 12099
$ws.Cells.Item(356, 3).Value = 12879
$ws.Cells.Item(356, 4).Value = 12086.4501953125
$ws.Cells.Item(356, 5).Value = 12505.9501953125
$ws.Cells.Item(356, 6).Value = 12505.9501953125
$ws.Cells.Item(356, 7).Value = 2339311
$ws.Cells.Item(356, 8).Value = 2024
$ws.Cells.Item(356, 9).Value = 7
$ws.Cells.Item(356, 10).Value = 1
$ws.Cells.Item(356, 11).Value = 0
$ws.Cells.Item(356, 12).Value = 0
$ws.Cells.Item(356, 13).Value = 0
$ws.Cells.Item(356, 14).Value = 27
$ws.Cells.Item(356, 15).Value = 1
$ws.Cells.Item(356, 16).Value = 0
$ws.Cells.Item(356, 17).Value = 0
$ws.Cells.Item(356, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# row 357
$ws.Cells.Item(357, 1).Value = 45481
$ws.Cells.Item(357, 2).Value = 12549
$ws.Cells.Item(357, 3).Value = 12689.400390625
$ws.Cells.Item(357, 4).Value = 12164.2998046875
$ws.Cells.Item(357, 5).Value = 12410
$ws.Cells.Item(357, 6).Value = 12410
$ws.Cells.Item(357, 7).Value = 1256082
$ws.Cells.Item(357, 8).Value = 2024
$ws.Cells.Item(357, 9).Value = 7
$ws.Cells.Item(357, 10).Value = 8
$ws.Cells.Item(357, 11).Value = 0
$ws.Cells.Item(357, 12).Value = 0
$ws.Cells.Item(357, 13).Value = 0
$ws.Cells.Item(357, 14).Value = 28
$ws.Cells.Item(357, 15).Value = 0
$ws.Cells.Item(357, 16).Value = 0
$ws.Cells.Item(357, 17).Value = 0
$ws.Cells.Item(357, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# row 358
$ws.Cells.Item(358, 1).Value = 45488
$ws.Cells.Item(358, 2).Value = 12500
$ws.Cells.Item(358, 3).Value = 12754.5498046875
$ws.Cells.Item(358, 4).Value = 11210.5498046875
$ws.Cells.Item(358, 5).Value = 11267.2001953125
$ws.Cells.Item(358, 6).Value = 11267.2001953125
$ws.Cells.Item(358, 7).Value = 1933896
$ws.Cells.Item(358, 8).Value = 2024
$ws.Cells.Item(358, 9).Value = 7
$ws.Cells.Item(358, 10).Value = 15
$ws.Cells.Item(358, 11).Value = 0
$ws.Cells.Item(358, 12).Value = 0
$ws.Cells.Item(358, 13).Value = 0
$ws.Cells.Item(358, 14).Value = 29
$ws.Cells.Item(358, 15).Value = 0
$ws.Cells.Item(358, 16).Value = 0
$ws.Cells.Item(358, 17).Value = 0
$ws.Cells.Item(358, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# row 359
$ws.Cells.Item(359, 1).Value = 45495
$ws.Cells.Item(359, 2).Value = 11267.2001953125
$ws.Cells.Item(359, 3).Value = 11475
$ws.Cells.Item(359, 4).Value = 10620
$ws.Cells.Item(359, 5).Value = 11272.599609375
$ws.Cells.Item(359, 6).Value = 11272.599609375
$ws.Cells.Item(359, 7).Value = 2285921
$ws.Cells.Item(359, 8).Value = 2024
$ws.Cells.Item(359, 9).Value = 7
$ws.Cells.Item(359, 10).Value = 22
$ws.Cells.Item(359, 11).Value = 0
$ws.Cells.Item(359, 12).Value = 0
$ws.Cells.Item(359, 13).Value = 0
$ws.Cells.Item(359, 14).Value = 30
$ws.Cells.Item(359, 15).Value = 0
$ws.Cells.Item(359, 16).Value = 0
$ws.Cells.Item(359, 17).Value = 0
$ws.Cells.Item(359, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# row 360
$ws.Cells.Item(360, 1).Value = 45502
$ws.Cells.Item(360, 2).Value = 11349.9501953125
$ws.Cells.Item(360, 3).Value = 12690
$ws.Cells.Item(360, 4).Value = 11225
$ws.Cells.Item(360, 5).Value = 11654.5498046875
$ws.Cells.Item(360, 6).Value = 11654.5498046875
$ws.Cells.Item(360, 7).Value = 4317411
$ws.Cells.Item(360, 8).Value = 2024
$ws.Cells.Item(360, 9).Value = 7
$ws.Cells.Item(360, 10).Value = 29
$ws.Cells.Item(360, 11).Value = 0
$ws.Cells.Item(360, 12).Value = 0
$ws.Cells.Item(360, 13).Value = 0
$ws.Cells.Item(360, 14).Value = 31
$ws.Cells.Item(360, 15).Value = 0
$ws.Cells.Item(360, 16).Value = 0
$ws.Cells.Item(360, 17).Value = 0
$ws.Cells.Item(360, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# row 361
$ws.Cells.Item(361, 1).Value = 45509
$ws.Cells.Item(361, 2).Value = 10950.0498046875
$ws.Cells.Item(361, 3).Value = 11800.0498046875
$ws.Cells.Item(361, 4).Value = 10950.0498046875
$ws.Cells.Item(361, 5).Value = 11740.4501953125
$ws.Cells.Item(361, 6).Value = 11740.4501953125
$ws.Cells.Item(361, 7).Value = 2166450
$ws.Cells.Item(361, 8).Value = 2024
$ws.Cells.Item(361, 9).Value = 8
$ws.Cells.Item(361, 10).Value = 5
$ws.Cells.Item(361, 11).Value = 0
$ws.Cells.Item(361, 12).Value = 0
$ws.Cells.Item(361, 13).Value = 0
$ws.Cells.Item(361, 14).Value = 32
$ws.Cells.Item(361, 15).Value = 0
$ws.Cells.Item(361, 16).Value = 0
$ws.Cells.Item(361, 17).Value = 0
$ws.Cells.Item(361, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

